# Update the "取得日時" (acquisition timestamp) column for rows 2-10
# on the "ランサーズ" sheet from 2025-12-19 12:37:51 to 2025-12-19 12:50:07,
# reflecting a later append run of the scraper.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-19 12:50:07"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
